# Updated cryptos list (refreshed prices / 1h volume %, and two rows whose
# coins swapped rank position). Numeric-looking Price values are written
# with a leading quote-prefix so Excel keeps them as literal text (matching
# the "67.799.71"-style dotted-thousands text already used in column D),
# then the cell style is reset to "Normal" so no stray number format sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.832.10'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '3.813.20'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''602.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("D6").Value = '''166.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '''0.519'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("D11").Value = '''6.37'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("E12").Value = '  -0.67%  '
$ws.Range("D13").Value = '''36.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").Value = '4.449.33'
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("D15").Value = '3.812.88'
$ws.Range("E15").Value = '  +1.29%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''18.49'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '67.843.04'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("E19").Value = '  +1.72%  '
$ws.Range("D20").Value = '''464.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("D21").Value = '''9.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D22").Value = '''0.704'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").Value = '''0.0000149'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.00%  '
$ws.Range("D24").Value = '''83.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("D25").Value = '''12.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").Value = '''10.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").Value = '3.962.67'
$ws.Range("E29").Value = '  +0.79%  '
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").Value = '''7.43'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.47%  '
$ws.Range("D32").Value = '''2.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("D33").Value = '''29.47'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").Value = '''9.09'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = '''0.1000'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("D38").Value = '''0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.37%  '
$ws.Range("D39").Value = '''5.82'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.99%  '
$ws.Range("D40").Value = '''3.25'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.39%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D43").Value = '''45.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").Value = '''47.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").Value = '''0.300'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").Value = '''28.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.17%  '
$ws.Range("D47").Value = '''151.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.74%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = '''1.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.35%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '''8.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("E50").Value = '  +2.12%  '
$ws.Range("D51").Value = '''392.42'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.07%  '
